# "bullet delete and balancing"
# PLAYER_BULLET_DISTANCE (row 2, col B): 30 -> 60
# PLAYER_BULLET_RELOAD   (row 4, col B): 60 -> 120
# Active cell selection moves from H17 to F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60
$ws.Range("B4").Value = 120

$ws.Range("F8").Select()
